$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "login" sheet: a handful of value tweaks (browser, url casing, locators)
# ---------------------------------------------------------------------------
$login = $wb.Worksheets.Item("login")

$login.Range("E2").Value = "firefox"
$login.Range("D3").Value = "enter url"
$login.Range("B4").Value = "name"
$login.Range("B5").Value = "id"
$login.Range("C7").Value = "Facebook"

# Column C got widened on this sheet too.
$login.Columns.Item(3).ColumnWidth = 36.2506

# ---------------------------------------------------------------------------
# "signup" sheet: the page got reworked into a full registration flow.
# Rows 1-3 (headers, open browser, launch url) are untouched; row 4 onward
# is replaced / extended with the new signup steps.
# ---------------------------------------------------------------------------
$signup = $wb.Worksheets.Item("signup")

$signup.Range("A4").Value = "Click On Create new account Button"
$signup.Range("B4").Value = "linkText"
$signup.Range("C4").Value = "Create New Account"
$signup.Range("D4").Value = "Click"
$signup.Range("E4").Value = "NA"

$signup.Range("A5").Value = "First Name"
$signup.Range("B5").Value = "xpath"
$signup.Range("C5").Value = "//input[@name='firstname']"
$signup.Range("D5").Value = "sendkeys"
$signup.Range("E5").Value = "Sinky"

$signup.Range("A6").Value = "SurName"
$signup.Range("B6").Value = "xpath"
$signup.Range("C6").Value = "//input[@name='lastname']"
$signup.Range("D6").Value = "sendkeys"
$signup.Range("E6").Value = "Verma"

$signup.Range("A7").Value = "MobileNumber or Email Address"
$signup.Range("B7").Value = "xpath"
$signup.Range("C7").Value = "//input[@name='reg_email__']"
$signup.Range("D7").Value = "sendkeys"
$signup.Range("E7").Value = "sinkykumariverma@gmail.com"

$signup.Range("A8").Value = "Re-Enter email address"
$signup.Range("B8").Value = "xpath"
$signup.Range("C8").Value = "//input[@name='reg_email_confirmation__']"
$signup.Range("D8").Value = "sendkeys"
$signup.Range("E8").Value = "sinkykumariverma@gmail.com"

$signup.Range("A9").Value = "New Password"
$signup.Range("B9").Value = "xpath"
$signup.Range("C9").Value = "//input[@name='reg_passwd__']"
$signup.Range("D9").Value = "sendkeys"
$signup.Range("E9").Value = "Sinky123"

$signup.Range("A10").Value = "Select date from drop down"
$signup.Range("B10").Value = "id"
$signup.Range("C10").Value = "day"
$signup.Range("D10").Value = "Click"
$signup.Range("E10").Value = 21

$signup.Range("A11").Value = "Select month from drop down"
$signup.Range("B11").Value = "id"
$signup.Range("C11").Value = "month"
$signup.Range("D11").Value = "Click"
$signup.Range("E11").Value = "April"

$signup.Range("A12").Value = "select year from drop down"
$signup.Range("B12").Value = "id"
$signup.Range("C12").Value = "year"
$signup.Range("D12").Value = "Click"
$signup.Range("E12").Value = 1999

$signup.Range("A13").Value = "Select Gender radio button"
$signup.Range("B13").Value = "xpath"
$signup.Range("C13").Value = "//label[text()='Female']"
$signup.Range("D13").Value = "Click"
$signup.Range("E13").Value = "Female"

$signup.Range("A14").Value = "Click on signUp button"
$signup.Range("B14").Value = "name"
$signup.Range("C14").Value = "websubmit"
$signup.Range("D14").Value = "Click"
$signup.Range("E14").Value = "NA"

$signup.Range("A15").Value = "Close Browser"
$signup.Range("B15").Value = "NA"
$signup.Range("C15").Value = "NA"
$signup.Range("D15").Value = "quit"
$signup.Range("E15").Value = "NA"

# New mailto hyperlinks for the two e-mail confirmation steps.
[void]$signup.Hyperlinks.Add($signup.Range("E7"), "mailto:sinkykumariverma@gmail.com")
$signup.Range("E7").Style = "Hyperlink"
[void]$signup.Hyperlinks.Add($signup.Range("E8"), "mailto:sinkykumariverma@gmail.com")
$signup.Range("E8").Style = "Hyperlink"

# Column widths re-tuned for the longer locator values / headers.
$signup.Columns.Item(1).ColumnWidth = 36.3756
$signup.Columns.Item(3).ColumnWidth = 49.4173
$signup.Columns.Item(4).ColumnWidth = 17.8056
$signup.Columns.Item(5).ColumnWidth = 27.9506

# ---------------------------------------------------------------------------
# Selection / active tab: the author left off with "signup" active, editing
# around E10, while "login" keeps a lingering selection on C7.
# ---------------------------------------------------------------------------
[void]$login.Range("C7").Select()
[void]$signup.Activate()
[void]$signup.Range("E10").Select()
